$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the first respondent's name
$ws.Range("A2").Value = "iiandjdmd"

# Remove the two duplicate "богдана" response rows (old rows 3 and 4);
# this shifts the summary rows below up by two (old 5/6/7 -> new 3/4/5).
$ws.Rows("3:4").Delete()

# Update the summary counts now that only one response row remains.
$ws.Range("B3").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 1
